$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting (values like "213.33"
# or "1.98" must not be auto-converted to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.621.78"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.649.39"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "213.33"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "0.535"
$ws.Range("E6").Value = "  +5.17%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "23.64"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "1.881.79"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.643.32"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "0.590"
$ws.Range("E14").Value = "  +5.36%  "
$ws.Range("D15").Value = "4.05"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "64.58"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "27.571.26"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "232.31"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "7.65"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "9.78"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "149.01"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "7.06"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "15.65"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").Value = "0.0487"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "3.20"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").Value = "1.432.67"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "0.573"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.888"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.818"
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("D43").Value = "5.50"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").Value = "65.17"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "1.791.36"
$ws.Range("D47").Value = "1.70"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "88.10"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -0.35%  "
